# Auto-generated script applying the 2024-10-22 crime-data update
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 6564
$ws.Range('K3').Value = 6774
$ws.Range('C4').Value = 1853
$ws.Range('F4').Value = 1919
$ws.Range('K4').Value = 1408
$ws.Range('K5').Value = 491
$ws.Range('K6').Value = 7443
$ws.Range('C7').Value = 28398
$ws.Range('F7').Value = 24112
$ws.Range('K7').Value = 22680

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K6').Value = 122
$ws.Range('K7').Value = 288

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 410
$ws.Range('K3').Value = 454
$ws.Range('K5').Value = 46
$ws.Range('K7').Value = 1489

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K3').Value = 173
$ws.Range('K5').Value = 11
$ws.Range('K7').Value = 489

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 251
$ws.Range('K3').Value = 352
$ws.Range('K6').Value = 308
$ws.Range('K7').Value = 989

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 220
$ws.Range('K6').Value = 229
$ws.Range('K7').Value = 770

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K3').Value = 129
$ws.Range('K6').Value = 193
$ws.Range('K7').Value = 532

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K3').Value = 158
$ws.Range('K7').Value = 380

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K6').Value = 161
$ws.Range('K7').Value = 679
$ws.Range('K8').Value = 1489
$ws.Range('K11').Value = 419
$ws.Range('K12').Value = 40
$ws.Range('K15').Value = 233
$ws.Range('K17').Value = 42
$ws.Range('K19').Value = 661
$ws.Range('K20').Value = 545
$ws.Range('K23').Value = 226
$ws.Range('K29').Value = 1230
$ws.Range('K31').Value = 250
$ws.Range('K33').Value = 989
$ws.Range('K36').Value = 289
$ws.Range('K37').Value = 770
$ws.Range('K42').Value = 837
$ws.Range('K44').Value = 187
$ws.Range('K46').Value = 46
$ws.Range('K47').Value = 151
$ws.Range('K49').Value = 124
$ws.Range('K51').Value = 286
$ws.Range('C52').Value = 649
$ws.Range('K52').Value = 602
$ws.Range('K53').Value = 288
$ws.Range('K54').Value = 447
$ws.Range('K55').Value = 244
$ws.Range('K57').Value = 84
$ws.Range('K60').Value = 133
$ws.Range('F63').Value = 203
$ws.Range('K63').Value = 60
$ws.Range('K64').Value = 141
$ws.Range('K65').Value = 532
$ws.Range('K67').Value = 888
$ws.Range('K73').Value = 204
$ws.Range('K74').Value = 24
$ws.Range('K78').Value = 257
$ws.Range('K79').Value = 566
$ws.Range('K80').Value = 81
$ws.Range('K83').Value = 489
$ws.Range('K84').Value = 184
$ws.Range('K85').Value = 1049
$ws.Range('K86').Value = 141
$ws.Range('K87').Value = 43
$ws.Range('K88').Value = 241
$ws.Range('K89').Value = 340
$ws.Range('K90').Value = 213
$ws.Range('K91').Value = 267
$ws.Range('K94').Value = 302
$ws.Range('K99').Value = 380
$ws.Range('C101').Value = 28398
$ws.Range('F101').Value = 24112
$ws.Range('K101').Value = 22680

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K6').Value = 86
$ws.Range('K7').Value = 250

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K3').Value = 325
$ws.Range('K7').Value = 888

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K3').Value = 75
$ws.Range('K7').Value = 184

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K3').Value = 26
$ws.Range('K7').Value = 124

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K3').Value = 107
$ws.Range('K7').Value = 447

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 350
$ws.Range('K3').Value = 434
$ws.Range('K6').Value = 359
$ws.Range('K7').Value = 1230

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 196
$ws.Range('K6').Value = 217
$ws.Range('K7').Value = 661

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K2').Value = 52
$ws.Range('K7').Value = 187

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('K6').Value = 44
$ws.Range('K7').Value = 161

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K3').Value = 257
$ws.Range('K7').Value = 837

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K3').Value = 64
$ws.Range('K4').Value = 23
$ws.Range('K7').Value = 257

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K2').Value = 74
$ws.Range('K7').Value = 244

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('K2').Value = 20
$ws.Range('K7').Value = 46

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K6').Value = 60
$ws.Range('K7').Value = 226

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K3').Value = 126
$ws.Range('K7').Value = 267

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K4').Value = 35
$ws.Range('K5').Value = 18
$ws.Range('K7').Value = 566

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K3').Value = 38
$ws.Range('K7').Value = 141

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K3').Value = 174
$ws.Range('K6').Value = 149
$ws.Range('K7').Value = 545

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range('K6').Value = 12
$ws.Range('K7').Value = 42

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K2').Value = 112
$ws.Range('K7').Value = 289

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K5').Value = 28
$ws.Range('K6').Value = 185
$ws.Range('K7').Value = 679

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K3').Value = 64
$ws.Range('K4').Value = 23
$ws.Range('K7').Value = 302

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K3').Value = 43
$ws.Range('K7').Value = 151

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K3').Value = 59
$ws.Range('K7').Value = 233

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K6').Value = 140
$ws.Range('K7').Value = 419

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K2').Value = 70
$ws.Range('K7').Value = 204

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K2').Value = 62
$ws.Range('K5').Value = 4
$ws.Range('K7').Value = 241

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K6').Value = 101
$ws.Range('K7').Value = 340

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K3').Value = 24
$ws.Range('K4').Value = 61
$ws.Range('K7').Value = 141

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K2').Value = 77
$ws.Range('K6').Value = 54
$ws.Range('K7').Value = 213

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K2').Value = 80
$ws.Range('K7').Value = 286

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('K6').Value = 36
$ws.Range('K7').Value = 84

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('K3').Value = 40
$ws.Range('K7').Value = 133

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 344
$ws.Range('K3').Value = 365
$ws.Range('K7').Value = 1049

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('K3').Value = 17
$ws.Range('K7').Value = 81

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K2').Value = 161
$ws.Range('K3').Value = 171
$ws.Range('C4').Value = 42
$ws.Range('K4').Value = 32
$ws.Range('C7').Value = 649
$ws.Range('K7').Value = 602

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('K3').Value = 11
$ws.Range('K6').Value = 14
$ws.Range('K7').Value = 40

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('K2').Value = 6
$ws.Range('K7').Value = 43

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range('K6').Value = 14
$ws.Range('K7').Value = 24
